# Remove the trailing asterisk from the "CLINMISKIN GEL*" product name
# in column A, for all rows where it appears (rows 2-202).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "CLINMISKIN GEL*") {
        $cell.Value2 = "CLINMISKIN GEL"
    }
}
